$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.643

$ws.Range("A4").Value = -21.984
$ws.Range("C4").Value = -12.78
$ws.Range("E4").Value = 12.647

$ws.Range("C5").Value = -12.894

$ws.Range("A6").Value = -21.345
$ws.Range("C6").Value = -12.652

$ws.Range("A7").Value = -21.118

$ws.Range("A8").Value = -21.826
$ws.Range("C8").Value = -13.026

$ws.Range("E9").Value = 13.018

$ws.Range("E11").Value = 12.852

$ws.Range("E14").Value = 13.056

$ws.Range("A16").Value = -21.114
$ws.Range("C16").Value = -12.192

$ws.Range("E18").Value = 12.657

$ws.Range("A20").Value = -21.948

$ws.Range("A21").Value = -20.921

$ws.Range("C22").Value = -12.78

$ws.Range("E25").Value = 12.464
